# ---------------------------------------------------------------------------
# feat: add 2022-Q1 data
#
# The source workbook rotates its "总计" (grand-total) sheet every quarter:
# the previous grand-total sheet is renamed to become the new quarter's
# per-fund holdings sheet, and a fresh "总计" sheet is appended at the end
# that re-summarises every quarter (including the newly added one).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# A sheet that already carries the canonical "fund holdings" header/row
# styling (bold, boxed header in row 1, boxed index column in col A) so we
# can clone that formatting instead of re-building it from scratch.
$styleTemplate = $wb.Worksheets.Item("2021-Q1")

# ---------------------------------------------------------------------------
# Step 1: repurpose the existing "总计" sheet into the new "2022-Q1" sheet.
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

$styleTemplate.Range("A2").Copy()
$q1.Range("A2:A11").PasteSpecial(-4122)
$styleTemplate.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$q1.Cells.Item(1,2).Value = "基金代码"
$q1.Cells.Item(1,3).Value = "基金名称"
$q1.Cells.Item(1,4).Value = "基金规模"
$q1.Cells.Item(1,5).Value = "股票总仓位"
$q1.Cells.Item(1,6).Value = "仓位占比"
$q1.Cells.Item(1,7).Value = "持有市值(亿元)"
$q1.Cells.Item(1,8).Value = "仓位排名"

# code, name, scale, stock position, position %, held market value (亿元), rank
$fundRows = @(
    ,@("960021","富兰克林国海潜力组合混合H 人民币","32.27","92.63","5.26","1.6974",5)
    ,@("450003","富兰克林国海潜力组合混合A 人民币","32.27","92.63","5.26","1.6974",5)
    ,@("450001","国富中国收益混合","22.23","61.42","4.23","0.9403",2)
    ,@("012510","富兰克林国海优质企业一年持有期混合型证券投资基金A","8.00","92.46","5.87","0.4696",1)
    ,@("010271","富兰克林国海价值成长一年持有期混合A","4.03","91.64","6.81","0.2744",1)
    ,@("010272","富兰克林国海价值成长一年持有期混合C","1.49","91.64","6.81","0.1015",1)
    ,@("400001","东方龙混合","2.80","84.04","3.61","0.1011",8)
    ,@("012511","富兰克林国海优质企业一年持有期混合型证券投资基金C","0.27","92.46","5.87","0.0158",1)
    ,@("014014","招商臻选平衡混合A","0.84","40.88","1.06","0.0089",8)
    ,@("014015","招商臻选平衡混合C","0.30","40.88","1.06","0.0032",8)
)

$r = 2
foreach ($row in $fundRows) {
    $q1.Cells.Item($r, 1).Value = ($r - 2)
    # fund code / scale / stock position / position% / market value are
    # stored as literal text in the source data (leading zeros in codes,
    # fixed two-decimal strings, etc.) - the leading "'" forces Excel to
    # keep them as text instead of auto-coercing to numbers.
    $q1.Cells.Item($r, 2).Value = "'" + $row[0]
    $q1.Cells.Item($r, 3).Value = $row[1]
    $q1.Cells.Item($r, 4).Value = "'" + $row[2]
    $q1.Cells.Item($r, 5).Value = "'" + $row[3]
    $q1.Cells.Item($r, 6).Value = "'" + $row[4]
    $q1.Cells.Item($r, 7).Value = "'" + $row[5]
    $q1.Cells.Item($r, 8).Value = $row[6]
    $r++
}

# Drop the implicit "quote prefix" formatting that typing a leading "'"
# leaves behind, so the text cells end up on the plain default style
# (matching every other quarterly sheet's data rows).
$q1.Range("B2:B11").ClearFormats()
$q1.Range("D2:G11").ClearFormats()

# ---------------------------------------------------------------------------
# Step 2: append a brand-new "总计" sheet at the end, re-summarising every
#         quarter including the newly added 2022-Q1 row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$total.Name = "总计"

$styleTemplate.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)
$styleTemplate.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)

$total.Cells.Item(1,2).Value = "日期"
$total.Cells.Item(1,3).Value = "持有数量(只)"
$total.Cells.Item(1,4).Value = "持有市值(亿元)"

# quarter label, holding count, held market value (亿元)
$totalRows = @(
    ,@("2022-Q1",10,5.31)
    ,@("2021-Q4",1,0.27)
    ,@("2021-Q3",2,0.72)
    ,@("2021-Q2",2,0.71)
    ,@("2021-Q1",1,0.29)
    ,@("2020-Q4",2,0.51)
)

$r = 2
foreach ($row in $totalRows) {
    $total.Cells.Item($r, 1).Value = ($r - 2)
    $total.Cells.Item($r, 2).Value = $row[0]
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]
    $r++
}

Write-Output "2022-Q1 sheet added; 总计 refreshed"
